$d = $word.ActiveDocument

# Delete the range covering paragraphs 3 through the end of the document
# (the numbered list of group names), leaving the title paragraph and the
# blank paragraph that follows it intact.
$startPara = $d.Paragraphs.Item(3)
$startRange = $startPara.Range.Start

$endRange = $d.Content.End

$range = $d.Range($startRange, $endRange)
$range.Delete()
